$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20 (pushes the totals block and everything
# below it down by one row), matching the "3 empty rows" layout change.
$ws.Rows.Item(20).Insert() | Out-Null
$ws.Rows.Item(20).RowHeight = 20.1

# Update the selection to reflect the author's post-edit cursor position.
$ws.Range("C19").Select() | Out-Null

# Keep the sheet's print area in sync with the new used range.
$ws.PageSetup.PrintArea = '$A$1:$K$54'

Write-Output "done"
